$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Numeric value updates (no style/type change) ---
$ws.Range("M14").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -27.272727272727
$ws.Range("I16").Value = 141
$ws.Range("J16").Value = 166
$ws.Range("K16").Value = -15.060240963855
$ws.Range("L16").Value = -9.615384615384
$ws.Range("M16").Value = -43.824701195219
$ws.Range("N16").Value = -85.553278688524
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 54.545454545454
$ws.Range("I17").Value = 405
$ws.Range("J17").Value = 335
$ws.Range("K17").Value = 20.895522388059
$ws.Range("L17").Value = 19.822485207100
$ws.Range("M17").Value = 105.583756345178
$ws.Range("N17").Value = -41.304347826087
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -13.333333333333
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = 6.060606060606
$ws.Range("L18").Value = 6.060606060606
$ws.Range("M18").Value = -31.372549019607
$ws.Range("N18").Value = -88.486842105263
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 70
$ws.Range("I19").Value = 385
$ws.Range("J19").Value = 328
$ws.Range("K19").Value = 17.378048780487
$ws.Range("L19").Value = -11.085450346420
$ws.Range("M19").Value = -25.675675675675
$ws.Range("N19").Value = -39.655172413793
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = -6.542056074766
$ws.Range("L20").Value = -4.761904761904
$ws.Range("M20").Value = -4.761904761904
$ws.Range("N20").Value = -90.808823529411
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -5
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = 28.048780487804
$ws.Range("I21").Value = 1180
$ws.Range("J21").Value = 1055
$ws.Range("K21").Value = 11.848341232227
$ws.Range("L21").Value = 2.966841186736
$ws.Range("M21").Value = -5.448717948717
$ws.Range("N21").Value = -72.929571002523
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 44
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = 62.962962962963
$ws.Range("L22").Value = 109.52380952381
$ws.Range("M22").Value = 51.724137931034
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 42.857142857142
$ws.Range("I23").Value = 143
$ws.Range("J23").Value = 129
$ws.Range("K23").Value = 10.852713178294
$ws.Range("L23").Value = -7.741935483870
$ws.Range("M23").Value = 26.548672566371
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 10.526315789473
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 42.253521126760
$ws.Range("I24").Value = 948
$ws.Range("J24").Value = 922
$ws.Range("K24").Value = 2.819956616052
$ws.Range("L24").Value = 5.099778270509
$ws.Range("M24").Value = -14.977578475336
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -45.161290322580
$ws.Range("I25").Value = 222
$ws.Range("J25").Value = 262
$ws.Range("K25").Value = -15.267175572519
$ws.Range("L25").Value = -17.164179104477
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -20
$ws.Range("I26").Value = 717
$ws.Range("J26").Value = 726
$ws.Range("K26").Value = -1.239669421487
$ws.Range("L26").Value = 27.353463587921
$ws.Range("M26").Value = 46.326530612244
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("M29").Value = -54.545454545454
$ws.Range("M30").Value = -47.058823529411

# --- Cells that become text placeholders ("0" / "***.*") with style matching other text cells ---
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C18").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# --- Cells that change FROM text placeholders TO numeric values ---
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$excel.CutCopyMode = $false
